$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) for new columns I and J - copy formatting from the
# existing header cell (H1) so the new headers match the bold/centered/
# bordered header style, then overwrite with the new header text.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-12 for columns I (I0) and J (IF)
$values = @{
    2  = @(9, 9)
    3  = @(7, 7)
    4  = @(6, 6)
    5  = @(7, 7)
    6  = @(6, 7)
    7  = @(7, 8)
    8  = @(7, 7)
    9  = @(7, 7)
    10 = @(6, 6)
    11 = @(9, 9)
    12 = @(7, 7)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
